$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: dataset re-pulled with a new row order and 2 additional
# records. Columns A,B,C,E,F,G,H,I,N,Q,R are constant for every record in this
# market/category subset, so only Fecha (D), Volumen (J), Precio minimo/maximo/
# promedio (K/L/M), Origen (O) and Precio $/Kg (P) need to move per row.
$rowUpdates = @(
  @(2, 44214, 40, 25000, 25000, 25000, "Región de La Araucanía", 1000),
  @(3, 44383, 80, 17000, 17000, 17000, "Provincia de Limarí", 680),
  @(4, 44160, 40, 11500, 11500, 11500, "Región del Maule", 460),
  @(5, 44390, 80, 16000, 16000, 16000, "Provincia de Limarí", 640),
  @(6, 44425, 90, 18000, 18000, 18000, "Provincia de Limarí", 720),
  @(7, 44386, 80, 17000, 17000, 17000, "Provincia de Limarí", 680),
  @(8, 44166, 120, 12000, 12000, 12000, "Región Metropolitana", 480),
  @(9, 44162, 140, 13000, 13000, 13000, "Región del Maule", 520),
  @(10, 44169, 160, 13000, 14000, 13500, "Región de La Araucanía", 540),
  @(11, 44203, 50, 22000, 22000, 22000, "Región de La Araucanía", 880),
  @(12, 44418, 90, 18000, 18000, 18000, "Provincia de Limarí", 720),
  @(13, 44172, 40, 12000, 12000, 12000, "Región de La Araucanía", 480),
  @(14, 44174, 20, 12000, 12000, 12000, "Región de La Araucanía", 480),
  @(15, 44407, 80, 18000, 18000, 18000, "Provincia de Limarí", 720),
  @(16, 44159, 150, 11500, 11500, 11500, "Región del Maule", 460),
  @(17, 44392, 90, 16000, 16000, 16000, "Provincia de Limarí", 640),
  @(18, 44362, 60, 20000, 20000, 20000, "Provincia de Limarí", 800),
  @(19, 44365, 80, 20000, 20000, 20000, "Provincia de Limarí", 800),
  @(20, 44358, 60, 20000, 20000, 20000, "Región Metropolitana", 800),
  @(21, 44421, 80, 17000, 17000, 17000, "Provincia de Limarí", 680),
  @(22, 44215, 60, 25000, 25000, 25000, "Región Metropolitana", 1000),
  @(23, 44414, 80, 18000, 18000, 18000, "Provincia de Limarí", 720),
  @(24, 44432, 80, 18000, 18000, 18000, "Provincia de Limarí", 720),
  @(25, 44351, 30, 20000, 20000, 20000, "Región Metropolitana", 800),
  @(26, 44369, 70, 18000, 18000, 18000, "Provincia de Limarí", 720),
  @(27, 44379, 70, 17000, 17000, 17000, "Provincia de Limarí", 680),
  @(28, 44446, 90, 17000, 17000, 17000, "Provincia de Limarí", 680),
  @(29, 44411, 80, 18000, 18000, 18000, "Provincia de Limarí", 720),
  @(30, 44161, 50, 11500, 11500, 11500, "Región del Maule", 460),
  @(31, 44428, 80, 18500, 18500, 18500, "Provincia de Limarí", 740),
  @(32, 44442, 80, 18000, 18000, 18000, "Provincia de Limarí", 720),
  @(33, 44435, 170, 18000, 19000, 18529, "Provincia de Limarí", 741),
  @(34, 44376, 70, 17000, 17000, 17000, "Provincia de Limarí", 680)
)

foreach ($row in $rowUpdates) {
  $r = $row[0]
  $ws.Cells.Item($r, 4).Value = $row[1]   # D: Fecha
  $ws.Cells.Item($r, 10).Value = $row[2]  # J: Volumen
  $ws.Cells.Item($r, 11).Value = $row[3]  # K: Precio minimo
  $ws.Cells.Item($r, 12).Value = $row[4]  # L: Precio maximo
  $ws.Cells.Item($r, 13).Value = $row[5]  # M: Precio promedio ponderado
  $ws.Cells.Item($r, 15).Value = $row[6]  # O: Origen
  $ws.Cells.Item($r, 16).Value = $row[7]  # P: Precio $/Kg
}

# --- Two brand-new records appended at the end of the subset ---
$newRows = @(
  @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44167, 10, 100112026, "Haba", "Sin especificar", "Primera", 30, 12000, 12000, 12000, "`$/saco 25 kilos", "Región de La Araucanía", 480, 25, "Hortaliza"),
  @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44400, 10, 100112026, "Haba", "Sin especificar", "Primera", 80, 16500, 16500, 16500, "`$/saco 25 kilos", "Provincia de Limarí", 660, 25, "Hortaliza")
)

$newRowIndex = 35
foreach ($row in $newRows) {
  for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item($newRowIndex, $col).Value = $row[$col - 1]
  }
  $ws.Cells.Item($newRowIndex, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $newRowIndex++
}
